# Cost Centre Master.xlsx - add "IsActive" column and drop the stray
# font-style override that was left on A1 ("Name").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 ("Name") was carrying an unused/erroneous cell style (applyFont) -
# clear it back to the default/unstyled cell.
$ws.Range("A1").ClearFormats()

# New header for the IsActive column.
$ws.Range("C1").Value = "IsActive"

# Match the author's final selection/cursor position.
$ws.Range("C2").Select()
